# Update cryptocurrency price and volume(1h) figures per the latest
# data refresh (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.960.29"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "'3.100.41"
$ws.Range("E3").Value = "  +3.04%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'579.70"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "'172.86"
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'3.095.17"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").Value = "'0.523"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").Value = "'6.43"
$ws.Range("E10").Value = "  -4.44%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "'37.27"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "'3.611.27"
$ws.Range("E16").Value = "  +3.67%  "
$ws.Range("D17").Value = "'66.897.54"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "'7.20"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "'3.099.07"
$ws.Range("E19").Value = "  +4.23%  "
$ws.Range("D20").Value = "'16.31"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").Value = "'483.40"
$ws.Range("E21").Value = "  +6.15%  "
$ws.Range("D22").Value = "'0.717"
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("E23").Value = "  +1.27%  "
$ws.Range("D24").Value = "'84.13"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").Value = "'13.25"
$ws.Range("E25").Value = "  +5.18%  "
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("B27").Value = "'Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("B28").Value = "'RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.02"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("D29").Value = "'7.99"
$ws.Range("E29").Value = "  -4.96%  "
$ws.Range("E30").Value = "  -2.46%  "
$ws.Range("D31").Value = "'2.68"
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("D32").Value = "'28.82"
$ws.Range("E32").Value = "  +3.18%  "
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("D37").Value = "'0.990"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "'47.71"
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  +3.08%  "
$ws.Range("D40").Value = "'50.11"
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("D41").Value = "'0.316"
$ws.Range("E41").Value = "  +2.14%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("D45").Value = "'2.849.88"
$ws.Range("E45").Value = "  +5.25%  "
$ws.Range("D46").Value = "'0.0359"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").Value = "'383.69"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("D48").Value = "'135.91"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("D50").Value = "'25.07"
$ws.Range("E50").Value = "  +2.36%  "
$ws.Range("E51").Value = "  -0.03%  "
